# BankAPP-Employee User stories.xlsx - updated user story file.
#
# Changes applied:
#  1. Insert a new column B ("Epic") before the existing "User Story" column
#     (which shifts the existing story text from column B to column C).
#  2. Remove the "check my account statement" story row entirely.
#  3. Renumber the "SI" column so it stays sequential (1..23).
#  4. Append two new user stories at the bottom ("Track my item" and
#     "Able to modify my item list").
#  5. Restore the selection to the cell the author ended up on (C6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Epic" column ---------------------------------------
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Epic"

# --- 2. Remove the "check my account statement" row -------------------------
# In the original layout this was SI #21, which lives at sheet row 22
# (row 1 is the header, row N holds SI #(N-1)).
$ws.Rows("22:22").Delete()

# --- 3. Re-sequence the SI column (A2:A23) so the numbers stay 1..22 -------
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- 4. Append the two new user stories at rows 23 and 24 -------------------
$ws.Cells.Item(23, 1).Value = 23 - 1
$ws.Cells.Item(23, 3).Value = "AS A Requester I WANT TO Track my item SO THAT i can easily Know when i get my item."

$ws.Cells.Item(24, 1).Value = 24 - 1
$ws.Cells.Item(24, 3).Value = "AS A Requester I WANT TO Able to modify my item list SO THAT i can do my work."

# --- 5. Restore selection ----------------------------------------------------
$ws.Range("C6").Select()
